$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Driver") and add the "I-Truck" header.
$null = $ws.Columns("C:C").Insert()
$ws.Range("C1").Value = "I-Truck"

# Match the original section's column width as closely as this engine's
# pixel-quantized ColumnWidth setter allows (target raw width 20.1640625).
$ws.Columns("C:C").ColumnWidth = 19.33

# Re-apply the autofilter so it covers the newly inserted column (A1:AA2).
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:AA2").AutoFilter()

# Update the _FilterDatabase defined names to the new extended range.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$AA`$2"
    }
}

# Restore the active selection recorded in the saved workbook.
$null = $ws.Range("C5").Select()
